$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price + 1h volume change) per upstream diff.
# Column D cells that are purely numeric-looking strings (e.g. "5.25") must be
# forced to Text so Excel does not silently convert them to numbers (losing
# trailing zeros / the original text formatting), matching the source data which
# stores these as formatted display strings, not numeric values.

$ws.Range("D2").Value = "37.466.41"
$ws.Range("E2").Value = "  -0.97%  "
$ws.Range("D3").Value = "2.050.70"
$ws.Range("E3").Value = "  -1.86%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.64"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.26%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  -3.69%  "
$ws.Range("E9").Value = "  -2.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0806"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.81%  "
$ws.Range("E11").Value = "  -2.04%  "
$ws.Range("D12").Value = "2.352.67"
$ws.Range("E12").Value = "  -1.86%  "
$ws.Range("E13").Value = "  -4.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.62"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -3.03%  "
$ws.Range("E15").Value = "  -3.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.25"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.27%  "
$ws.Range("D17").Value = "2.069.99"
$ws.Range("E17").Value = "  -1.24%  "
$ws.Range("D18").Value = "37.358.42"
$ws.Range("E18").Value = "  -1.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.09"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.08%  "
$ws.Range("E20").Value = "  -1.75%  "
$ws.Range("D21").Value = "0.0₃0843"
$ws.Range("E21").Value = "  +0.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "225.66"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.83%  "
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("E24").Value = "  -0.95%  "
$ws.Range("E25").Value = "  -4.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.52"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.55%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "168.48"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.87%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.129"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -3.85%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.37"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "18.89"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.33%  "
$ws.Range("E31").Value = "  -2.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.54"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.35%  "
$ws.Range("E33").Value = "  -3.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.56"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.60%  "
$ws.Range("E35").Value = "  -3.93%  "
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.20"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.42"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.39%  "
$ws.Range("E40").Value = "  -5.83%  "
$ws.Range("D41").Value = "1.500.56"
$ws.Range("E41").Value = "  +3.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "16.91"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.63%  "
$ws.Range("E43").Value = "  -1.70%  "
$ws.Range("E44").Value = "  -5.25%  "
$ws.Range("E45").Value = "  -4.15%  "
$ws.Range("E46").Value = "  -4.28%  "
$ws.Range("E47").Value = "  -4.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.22"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.07%  "
$ws.Range("B49").Value = "MXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.92"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.63%  "
$ws.Range("B50").Value = "FTXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.79"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -7.84%  "
$ws.Range("D51").Value = "2.238.63"
$ws.Range("E51").Value = "  -1.89%  "
